# Edit script for Belvo_Data_Dictionary_es.xlsx
# 1) links sheet B12: wording tweak
# 2) financial_statements sheet: insert new rows + append ".01" to Example values

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# 1) links sheet: tweak credentials_storage description wording
# -------------------------------------------------------------------------
$wsLinks = $wb.Worksheets.Item("links")
$wsLinks.Cells.Item(12,2).Value2 = "Indica si se deben almacenar las credenciales (y la duración durante la cual se almacenarán las credenciales).`n`n- Para enlaces recurrentes, esto se establece en ``store`` por defecto (y no se puede cambiar).`n- Para enlaces únicos, esto se establece en ``365d`` por defecto.`n`nPuede ser:`n  - ``store`` para almacenar credenciales (hasta que se elimine el enlace)`n  - ``nostore`` para no almacenar credenciales`n  - Cualquier valor entre ``1d`` y ``365d`` para indicar el número de días que deseas que se almacenen las credenciales.`n`nPara más información, consulta la sección <a href=`"https://developers.belvo.com/docs/data-retention-controls#credentials_storage`" target=`"_blank`">credentials_storage</a> de nuestro artículo sobre controles de retención de datos."

# -------------------------------------------------------------------------
# 2) financial_statements sheet
# -------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("financial_statements")

# --- Insert the 3 blocks of new rows (order matters: top-to-bottom) ---
$ws.Rows.Item(24).Insert()
$ws.Rows.Item(53).Resize(6).Insert()
$ws.Rows.Item(71).Resize(2).Insert()

# --- Populate newly inserted rows with full content ---
# row 24: balance_sheet.non_current_assets.accumulated_depreciation_and_amortization
$ws.Cells.Item(24,1).Value2 = "balance_sheet.non_current_assets.accumulated_depreciation_and_amortization"
$ws.Cells.Item(24,2).Value2 = "La depreciación y amortización acumulada total, que representa la asignación acumulativa del costo de los activos no corrientes durante el período en que se espera que proporcionen beneficios económicos."
$ws.Cells.Item(24,3).Value2 = "123456.01"
$ws.Cells.Item(24,4).Value2 = "number"
$ws.Cells.Item(24,5).Value2 = "float"
$ws.Cells.Item(24,6).Value2 = "Yes"
$ws.Cells.Item(24,7).Value2 = "Yes"

# row 54: balance_sheet.equity.future_capital_contributions
$ws.Cells.Item(54,1).Value2 = "balance_sheet.equity.future_capital_contributions"
$ws.Cells.Item(54,2).Value2 = "Los fondos recibidos de los accionistas que están específicamente designados para futuros aumentos de capital o inversiones."
$ws.Cells.Item(54,3).Value2 = "75000.01"
$ws.Cells.Item(54,4).Value2 = "number"
$ws.Cells.Item(54,5).Value2 = "float"
$ws.Cells.Item(54,6).Value2 = "Yes"
$ws.Cells.Item(54,7).Value2 = "Yes"

# row 55: balance_sheet.equity.legal_reserve
$ws.Cells.Item(55,1).Value2 = "balance_sheet.equity.legal_reserve"
$ws.Cells.Item(55,2).Value2 = "La reserva legal exigida por la ley, generalmente apartada de las ganancias, para proporcionar protección financiera contra pérdidas u obligaciones futuras."
$ws.Cells.Item(55,3).Value2 = "25000.01"
$ws.Cells.Item(55,4).Value2 = "number"
$ws.Cells.Item(55,5).Value2 = "float"
$ws.Cells.Item(55,6).Value2 = "Yes"
$ws.Cells.Item(55,7).Value2 = "Yes"

# row 56: balance_sheet.equity.capital_update_excess
$ws.Cells.Item(56,1).Value2 = "balance_sheet.equity.capital_update_excess"
$ws.Cells.Item(56,2).Value2 = "El excedente resultante de los ajustes realizados al capital social, a menudo debido a la inflación o la revalorización de activos."
$ws.Cells.Item(56,3).Value2 = "15000.01"
$ws.Cells.Item(56,4).Value2 = "number"
$ws.Cells.Item(56,5).Value2 = "float"
$ws.Cells.Item(56,6).Value2 = "Yes"
$ws.Cells.Item(56,7).Value2 = "Yes"

# row 57: balance_sheet.equity.capital_update_insufficiency
$ws.Cells.Item(57,1).Value2 = "balance_sheet.equity.capital_update_insufficiency"
$ws.Cells.Item(57,2).Value2 = "El déficit resultante de los ajustes realizados al capital social, a menudo debido a la inflación o la revalorización de activos."
$ws.Cells.Item(57,3).Value2 = "-5000.01"
$ws.Cells.Item(57,4).Value2 = "number"
$ws.Cells.Item(57,5).Value2 = "float"
$ws.Cells.Item(57,6).Value2 = "Yes"
$ws.Cells.Item(57,7).Value2 = "Yes"

# row 58: balance_sheet.equity.capital_reserve
$ws.Cells.Item(58,1).Value2 = "balance_sheet.equity.capital_reserve"
$ws.Cells.Item(58,2).Value2 = "La reserva de capital derivada de actividades no operativas, como ganancias de revalorizaciones de activos o ciertas transacciones de capital."
$ws.Cells.Item(58,3).Value2 = "10000.01"
$ws.Cells.Item(58,4).Value2 = "number"
$ws.Cells.Item(58,5).Value2 = "float"
$ws.Cells.Item(58,6).Value2 = "Yes"
$ws.Cells.Item(58,7).Value2 = "Yes"

# row 59: balance_sheet.equity.share_premium_on_stock_sales
$ws.Cells.Item(59,1).Value2 = "balance_sheet.equity.share_premium_on_stock_sales"
$ws.Cells.Item(59,2).Value2 = "El monto excedente recibido por una empresa cuando las acciones se emiten a un precio superior a su valor nominal (par)."
$ws.Cells.Item(59,3).Value2 = "50000.01"
$ws.Cells.Item(59,4).Value2 = "number"
$ws.Cells.Item(59,5).Value2 = "float"
$ws.Cells.Item(59,6).Value2 = "Yes"
$ws.Cells.Item(59,7).Value2 = "Yes"

# row 78: income_statement.income_statement_financial_gains
$ws.Cells.Item(78,1).Value2 = "income_statement.income_statement_financial_gains"
$ws.Cells.Item(78,2).Value2 = "El ingreso financiero total positivo, incluyendo ingresos por intereses, ganancias por diferencias de cambio y otras ganancias de actividades de financiamiento. Este valor siempre debe ser positivo."
$ws.Cells.Item(78,3).Value2 = "85000.01"
$ws.Cells.Item(78,4).Value2 = "number"
$ws.Cells.Item(78,5).Value2 = "float"
$ws.Cells.Item(78,6).Value2 = "Yes"
$ws.Cells.Item(78,7).Value2 = "Yes"

# row 79: income_statement.income_statement_financial_costs
$ws.Cells.Item(79,1).Value2 = "income_statement.income_statement_financial_costs"
$ws.Cells.Item(79,2).Value2 = "Los gastos financieros totales, incluidos los gastos por intereses, las pérdidas por diferencias de cambio y otros costos incurridos por actividades de financiación. Este valor siempre debe ser negativo."
$ws.Cells.Item(79,3).Value2 = "-32000.01"
$ws.Cells.Item(79,4).Value2 = "number"
$ws.Cells.Item(79,5).Value2 = "float"
$ws.Cells.Item(79,6).Value2 = "Yes"
$ws.Cells.Item(79,7).Value2 = "Yes"

# --- Append ".01" suffix to Example (column C) for all pre-existing rows ---
# (rows 11-23 keep their original row numbers; rows >=25 in final layout
#  correspond to pre-existing rows that have shifted down due to the inserts above)
$ws.Cells.Item(11,3).Value2 = "48572.01"
$ws.Cells.Item(12,3).Value2 = "21345.01"
$ws.Cells.Item(13,3).Value2 = "154321.01"
$ws.Cells.Item(14,3).Value2 = "31789.01"
$ws.Cells.Item(15,3).Value2 = "12345.01"
$ws.Cells.Item(16,3).Value2 = "0.01"
$ws.Cells.Item(17,3).Value2 = "8976.01"
$ws.Cells.Item(18,3).Value2 = "65432.01"
$ws.Cells.Item(19,3).Value2 = "14321.01"
$ws.Cells.Item(20,3).Value2 = "54321.01"
$ws.Cells.Item(21,3).Value2 = "372480.01"
$ws.Cells.Item(23,3).Value2 = "1123456.01"
$ws.Cells.Item(25,3).Value2 = "10987.01"
$ws.Cells.Item(26,3).Value2 = "5432.01"
$ws.Cells.Item(27,3).Value2 = "47654.01"
$ws.Cells.Item(28,3).Value2 = "43210.01"
$ws.Cells.Item(29,3).Value2 = "65432.01"
$ws.Cells.Item(30,3).Value2 = "32876.01"
$ws.Cells.Item(31,3).Value2 = "1346647.01"
$ws.Cells.Item(33,3).Value2 = "49876.01"
$ws.Cells.Item(34,3).Value2 = "103298.01"
$ws.Cells.Item(35,3).Value2 = "25643.01"
$ws.Cells.Item(36,3).Value2 = "14321.01"
$ws.Cells.Item(37,3).Value2 = "21987.01"
$ws.Cells.Item(38,3).Value2 = "12765.01"
$ws.Cells.Item(39,3).Value2 = "18765.01"
$ws.Cells.Item(40,3).Value2 = "10987.01"
$ws.Cells.Item(41,3).Value2 = "5321.01"
$ws.Cells.Item(42,3).Value2 = "260963.01"
$ws.Cells.Item(44,3).Value2 = "30876.01"
$ws.Cells.Item(45,3).Value2 = "42310.01"
$ws.Cells.Item(46,3).Value2 = "21987.01"
$ws.Cells.Item(47,3).Value2 = "10987.01"
$ws.Cells.Item(48,3).Value2 = "26543.01"
$ws.Cells.Item(49,3).Value2 = "30218.01"
$ws.Cells.Item(50,3).Value2 = "15432.01"
$ws.Cells.Item(51,3).Value2 = "178353.01"
$ws.Cells.Item(53,3).Value2 = "501234.01"
$ws.Cells.Item(60,3).Value2 = "202345.01"
$ws.Cells.Item(61,3).Value2 = "10987.01"
$ws.Cells.Item(62,3).Value2 = "70876.01"
$ws.Cells.Item(63,3).Value2 = "50321.01"
$ws.Cells.Item(64,3).Value2 = "836763.01"
$ws.Cells.Item(66,3).Value2 = "1212345.01"
$ws.Cells.Item(67,3).Value2 = "1123456.01"
$ws.Cells.Item(68,3).Value2 = "88987.01"
$ws.Cells.Item(69,3).Value2 = "609876.01"
$ws.Cells.Item(70,3).Value2 = "412345.01"
$ws.Cells.Item(71,3).Value2 = "101234.01"
$ws.Cells.Item(72,3).Value2 = "190890.01"
$ws.Cells.Item(74,3).Value2 = "122345.01"
$ws.Cells.Item(75,3).Value2 = "68545.01"
$ws.Cells.Item(77,3).Value2 = "15098.01"
$ws.Cells.Item(80,3).Value2 = "5678.01"
$ws.Cells.Item(81,3).Value2 = "89321.01"
$ws.Cells.Item(83,3).Value2 = "20123.01"
$ws.Cells.Item(84,3).Value2 = "69198.01"
$ws.Cells.Item(86,3).Value2 = "0.01"
$ws.Cells.Item(87,3).Value2 = "69198.01"
